$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared "Commerce" section label (was "Commerce ; reparation...") ---
$ws.Range("H8").Value = "Commerce"
$ws.Range("H27").Value = "Commerce"
$ws.Range("H46").Value = "Commerce"
$ws.Range("H63").Value = "Commerce"
$ws.Range("H79").Value = "Commerce"
$ws.Range("H97").Value = "Commerce"
$ws.Range("H114").Value = "Commerce"
$ws.Range("H130").Value = "Commerce"
$ws.Range("H147").Value = "Commerce"
$ws.Range("H164").Value = "Commerce"
$ws.Range("H180").Value = "Commerce"
$ws.Range("H196").Value = "Commerce"
$ws.Range("H211").Value = "Commerce"
$ws.Range("H229").Value = "Commerce"
$ws.Range("H247").Value = "Commerce"
$ws.Range("H266").Value = "Commerce"
$ws.Range("H284").Value = "Commerce"
$ws.Range("H301").Value = "Commerce"

# --- Update nombre_aides / montant_total for the refreshed data pull ---
$ws.Range("C2").Value = 3278
$ws.Range("D2").Value = 4441450
$ws.Range("C4").Value = 7587
$ws.Range("D4").Value = 9608658
$ws.Range("C5").Value = 40
$ws.Range("D5").Value = 53995
$ws.Range("C6").Value = 372
$ws.Range("D6").Value = 461954
$ws.Range("C7").Value = 21367
$ws.Range("D7").Value = 30811050
$ws.Range("C8").Value = 28583
$ws.Range("D8").Value = 38348688
$ws.Range("C9").Value = 5818
$ws.Range("D9").Value = 8080271
$ws.Range("C10").Value = 22278
$ws.Range("D10").Value = 32074313
$ws.Range("C11").Value = 2889
$ws.Range("D11").Value = 3787019
$ws.Range("C12").Value = 1327
$ws.Range("D12").Value = 1873320
$ws.Range("C13").Value = 4940
$ws.Range("D13").Value = 6800742
$ws.Range("C14").Value = 15278
$ws.Range("D14").Value = 20676082
$ws.Range("C15").Value = 7350
$ws.Range("D15").Value = 9499799
$ws.Range("C17").Value = 16576
$ws.Range("D17").Value = 21249746
$ws.Range("C18").Value = 18937
$ws.Range("D18").Value = 25473512
$ws.Range("C19").Value = 7982
$ws.Range("D19").Value = 9760161
$ws.Range("C20").Value = 24436
$ws.Range("D20").Value = 30072010
$ws.Range("C21").Value = 1207
$ws.Range("D21").Value = 1642208
$ws.Range("C23").Value = 2142
$ws.Range("D23").Value = 2627953
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 6000
$ws.Range("C25").Value = 166
$ws.Range("D25").Value = 208030
$ws.Range("C26").Value = 5294
$ws.Range("D26").Value = 7599416
$ws.Range("C27").Value = 9180
$ws.Range("D27").Value = 12221314
$ws.Range("C28").Value = 776
$ws.Range("D28").Value = 1050834
$ws.Range("C29").Value = 6633
$ws.Range("D29").Value = 9459020
$ws.Range("C30").Value = 558
$ws.Range("D30").Value = 693786
$ws.Range("C31").Value = 348
$ws.Range("D31").Value = 469871
$ws.Range("C32").Value = 1126
$ws.Range("D32").Value = 1505837
$ws.Range("C33").Value = 3006
$ws.Range("D33").Value = 4053400
$ws.Range("C34").Value = 1621
$ws.Range("D34").Value = 2045895
$ws.Range("C36").Value = 2346
$ws.Range("D36").Value = 2816925
$ws.Range("C37").Value = 4160
$ws.Range("D37").Value = 5571047
$ws.Range("C38").Value = 1954
$ws.Range("D38").Value = 2320802
$ws.Range("C39").Value = 7713
$ws.Range("D39").Value = 9631322
$ws.Range("C40").Value = 1096
$ws.Range("D40").Value = 1495061
$ws.Range("C42").Value = 2811
$ws.Range("D42").Value = 3365143
$ws.Range("C43").Value = 7
$ws.Range("D43").Value = 10342
$ws.Range("C44").Value = 116
$ws.Range("D44").Value = 138201
$ws.Range("C45").Value = 5911
$ws.Range("D45").Value = 8478948
$ws.Range("C46").Value = 10058
$ws.Range("D46").Value = 13544853
$ws.Range("C47").Value = 1036
$ws.Range("D47").Value = 1394176
$ws.Range("C48").Value = 8850
$ws.Range("D48").Value = 12700210
$ws.Range("C49").Value = 854
$ws.Range("D49").Value = 1070575
$ws.Range("C50").Value = 394
$ws.Range("D50").Value = 555951
$ws.Range("C51").Value = 1561
$ws.Range("D51").Value = 2055551
$ws.Range("C52").Value = 4461
$ws.Range("D52").Value = 5997926
$ws.Range("C53").Value = 2002
$ws.Range("D53").Value = 2487581
$ws.Range("C54").Value = 2535
$ws.Range("D54").Value = 3100403
$ws.Range("C55").Value = 6758
$ws.Range("D55").Value = 9173306
$ws.Range("C56").Value = 2400
$ws.Range("D56").Value = 2846497
$ws.Range("C57").Value = 8551
$ws.Range("D57").Value = 10704206
$ws.Range("C58").Value = 853
$ws.Range("D58").Value = 1166613
$ws.Range("C59").Value = 1751
$ws.Range("D59").Value = 2203280
$ws.Range("C61").Value = 95
$ws.Range("D61").Value = 99170
$ws.Range("C62").Value = 5274
$ws.Range("D62").Value = 7533339
$ws.Range("C63").Value = 7487
$ws.Range("D63").Value = 9829085
$ws.Range("C64").Value = 939
$ws.Range("D64").Value = 1320317
$ws.Range("C65").Value = 5239
$ws.Range("D65").Value = 7495893
$ws.Range("C66").Value = 577
$ws.Range("D66").Value = 740402
$ws.Range("C67").Value = 294
$ws.Range("D67").Value = 417445
$ws.Range("C68").Value = 1167
$ws.Range("D68").Value = 1574723
$ws.Range("C69").Value = 2854
$ws.Range("D69").Value = 3832721
$ws.Range("C70").Value = 1546
$ws.Range("D70").Value = 1945204
$ws.Range("C71").Value = 1877
$ws.Range("D71").Value = 2338673
$ws.Range("C72").Value = 3529
$ws.Range("D72").Value = 4699647
$ws.Range("C73").Value = 1804
$ws.Range("D73").Value = 2179349
$ws.Range("C74").Value = 6939
$ws.Range("D74").Value = 8651126
$ws.Range("C75").Value = 766
$ws.Range("D75").Value = 1055323
$ws.Range("C76").Value = 686
$ws.Range("D76").Value = 902916
$ws.Range("C78").Value = 1704
$ws.Range("D78").Value = 2427054
$ws.Range("C79").Value = 1772
$ws.Range("D79").Value = 2445768
$ws.Range("C80").Value = 235
$ws.Range("D80").Value = 330520
$ws.Range("C81").Value = 1581
$ws.Range("D81").Value = 2256448
$ws.Range("C82").Value = 118
$ws.Range("D82").Value = 154938
$ws.Range("C83").Value = 46
$ws.Range("D83").Value = 63418
$ws.Range("C84").Value = 352
$ws.Range("D84").Value = 477088
$ws.Range("C85").Value = 1166
$ws.Range("D85").Value = 1625273
$ws.Range("C86").Value = 690
$ws.Range("D86").Value = 871863
$ws.Range("C87").Value = 366
$ws.Range("D87").Value = 447584
$ws.Range("C88").Value = 681
$ws.Range("D88").Value = 908027
$ws.Range("C89").Value = 344
$ws.Range("D89").Value = 435737
$ws.Range("C90").Value = 1248
$ws.Range("D90").Value = 1474749
$ws.Range("C91").Value = 1549
$ws.Range("D91").Value = 2128555
$ws.Range("C93").Value = 3620
$ws.Range("D93").Value = 4639536
$ws.Range("C94").Value = 8
$ws.Range("D94").Value = 10810
$ws.Range("C95").Value = 283
$ws.Range("D95").Value = 330720
$ws.Range("C96").Value = 10931
$ws.Range("D96").Value = 15769669
$ws.Range("C97").Value = 16182
$ws.Range("D97").Value = 21697141
$ws.Range("C98").Value = 2117
$ws.Range("D98").Value = 2900319
$ws.Range("C99").Value = 10679
$ws.Range("D99").Value = 15444833
$ws.Range("C100").Value = 1220
$ws.Range("D100").Value = 1558243
$ws.Range("C101").Value = 693
$ws.Range("D101").Value = 968654
$ws.Range("C102").Value = 1838
$ws.Range("D102").Value = 2509642
$ws.Range("C103").Value = 6468
$ws.Range("D103").Value = 8682368
$ws.Range("C104").Value = 3051
$ws.Range("D104").Value = 3868592
$ws.Range("C105").Value = 3868
$ws.Range("D105").Value = 4725763
$ws.Range("C106").Value = 7901
$ws.Range("D106").Value = 10790752
$ws.Range("C107").Value = 3368
$ws.Range("D107").Value = 3977771
$ws.Range("C108").Value = 15481
$ws.Range("D108").Value = 18704620
$ws.Range("C109").Value = 637
$ws.Range("D109").Value = 888950
$ws.Range("C110").Value = 794
$ws.Range("D110").Value = 1106715
$ws.Range("C112").Value = 19
$ws.Range("D112").Value = 28000
$ws.Range("C113").Value = 1343
$ws.Range("D113").Value = 1937450
$ws.Range("C114").Value = 2194
$ws.Range("D114").Value = 3036964
$ws.Range("C115").Value = 533
$ws.Range("D115").Value = 757960
$ws.Range("C116").Value = 1316
$ws.Range("D116").Value = 1881576
$ws.Range("C117").Value = 214
$ws.Range("D117").Value = 288915
$ws.Range("C118").Value = 70
$ws.Range("D118").Value = 97706
$ws.Range("C119").Value = 302
$ws.Range("D119").Value = 420038
$ws.Range("C120").Value = 872
$ws.Range("D120").Value = 1209946
$ws.Range("C121").Value = 918
$ws.Range("D121").Value = 1226268
$ws.Range("C122").Value = 498
$ws.Range("D122").Value = 634691
$ws.Range("C123").Value = 753
$ws.Range("D123").Value = 1057231
$ws.Range("C124").Value = 302
$ws.Range("D124").Value = 397547
$ws.Range("C125").Value = 925
$ws.Range("D125").Value = 1170932
$ws.Range("C126").Value = 280
$ws.Range("D126").Value = 383126
$ws.Range("C127").Value = 144
$ws.Range("D127").Value = 201952
$ws.Range("C129").Value = 299
$ws.Range("D129").Value = 437270
$ws.Range("C130").Value = 403
$ws.Range("D130").Value = 572026
$ws.Range("C131").Value = 96
$ws.Range("D131").Value = 140079
$ws.Range("C132").Value = 270
$ws.Range("D132").Value = 384237
$ws.Range("C133").Value = 38
$ws.Range("D133").Value = 54360
$ws.Range("C135").Value = 25
$ws.Range("D135").Value = 35905
$ws.Range("C136").Value = 193
$ws.Range("D136").Value = 273522
$ws.Range("C137").Value = 155
$ws.Range("D137").Value = 213519
$ws.Range("C138").Value = 107
$ws.Range("D138").Value = 150336
$ws.Range("C139").Value = 124
$ws.Range("D139").Value = 175400
$ws.Range("C140").Value = 54
$ws.Range("D140").Value = 74194
$ws.Range("C141").Value = 162
$ws.Range("D141").Value = 208857
$ws.Range("C142").Value = 976
$ws.Range("D142").Value = 1302031
$ws.Range("C143").Value = 2950
$ws.Range("D143").Value = 3593804
$ws.Range("C144").Value = 11
$ws.Range("D144").Value = 15525
$ws.Range("C145").Value = 153
$ws.Range("D145").Value = 181960
$ws.Range("C146").Value = 9604
$ws.Range("D146").Value = 13614980
$ws.Range("C147").Value = 16536
$ws.Range("D147").Value = 21228444
$ws.Range("C148").Value = 3512
$ws.Range("D148").Value = 4843812
$ws.Range("C149").Value = 10512
$ws.Range("D149").Value = 15194870
$ws.Range("C150").Value = 1133
$ws.Range("D150").Value = 1448287
$ws.Range("C151").Value = 723
$ws.Range("D151").Value = 1015117
$ws.Range("C152").Value = 1876
$ws.Range("D152").Value = 2528971
$ws.Range("C153").Value = 6640
$ws.Range("D153").Value = 8971277
$ws.Range("C154").Value = 3172
$ws.Range("D154").Value = 3982009
$ws.Range("C155").Value = 3532
$ws.Range("D155").Value = 4400370
$ws.Range("C156").Value = 8061
$ws.Range("D156").Value = 10960531
$ws.Range("C157").Value = 3371
$ws.Range("D157").Value = 4028183
$ws.Range("C158").Value = 15414
$ws.Range("D158").Value = 18087035
$ws.Range("C159").Value = 892
$ws.Range("D159").Value = 1201596
$ws.Range("C160").Value = 1444
$ws.Range("D160").Value = 1934107
$ws.Range("C162").Value = 30
$ws.Range("D162").Value = 44673
$ws.Range("C163").Value = 2763
$ws.Range("D163").Value = 4030612
$ws.Range("C164").Value = 4892
$ws.Range("D164").Value = 6713915
$ws.Range("C165").Value = 628
$ws.Range("D165").Value = 884867
$ws.Range("C166").Value = 2425
$ws.Range("D166").Value = 3423100
$ws.Range("C167").Value = 317
$ws.Range("D167").Value = 430040
$ws.Range("C168").Value = 137
$ws.Range("D168").Value = 198483
$ws.Range("C169").Value = 495
$ws.Range("D169").Value = 696255
$ws.Range("C170").Value = 1474
$ws.Range("D170").Value = 2086537
$ws.Range("C171").Value = 839
$ws.Range("D171").Value = 1160612
$ws.Range("C172").Value = 1320
$ws.Range("D172").Value = 1711970
$ws.Range("C173").Value = 1784
$ws.Range("D173").Value = 2518928
$ws.Range("C174").Value = 714
$ws.Range("D174").Value = 934892
$ws.Range("C175").Value = 2168
$ws.Range("D175").Value = 2774319
$ws.Range("C176").Value = 367
$ws.Range("D176").Value = 469058
$ws.Range("C177").Value = 602
$ws.Range("D177").Value = 805284
$ws.Range("C179").Value = 839
$ws.Range("D179").Value = 1210877
$ws.Range("C180").Value = 1603
$ws.Range("D180").Value = 2208655
$ws.Range("C181").Value = 488
$ws.Range("D181").Value = 681057
$ws.Range("C182").Value = 872
$ws.Range("D182").Value = 1260052
$ws.Range("C183").Value = 143
$ws.Range("D183").Value = 187381
$ws.Range("C184").Value = 62
$ws.Range("D184").Value = 90896
$ws.Range("C185").Value = 190
$ws.Range("D185").Value = 269726
$ws.Range("C186").Value = 766
$ws.Range("D186").Value = 1066565
$ws.Range("C187").Value = 675
$ws.Range("D187").Value = 919425
$ws.Range("C188").Value = 414
$ws.Range("D188").Value = 533888
$ws.Range("C189").Value = 676
$ws.Range("D189").Value = 936239
$ws.Range("C190").Value = 258
$ws.Range("D190").Value = 336378
$ws.Range("C191").Value = 942
$ws.Range("D191").Value = 1221159
$ws.Range("C192").Value = 42
$ws.Range("D192").Value = 59714
$ws.Range("C193").Value = 61
$ws.Range("D193").Value = 86624
$ws.Range("C195").Value = 154
$ws.Range("D195").Value = 225741
$ws.Range("C196").Value = 396
$ws.Range("D196").Value = 556861
$ws.Range("C197").Value = 105
$ws.Range("D197").Value = 154459
$ws.Range("C198").Value = 126
$ws.Range("D198").Value = 183941
$ws.Range("C199").Value = 18
$ws.Range("D199").Value = 25043
$ws.Range("C201").Value = 63
$ws.Range("D201").Value = 94086
$ws.Range("C202").Value = 46
$ws.Range("D202").Value = 66893
$ws.Range("C203").Value = 62
$ws.Range("D203").Value = 92725
$ws.Range("C204").Value = 35
$ws.Range("D204").Value = 50992
$ws.Range("C205").Value = 30
$ws.Range("D205").Value = 43200
$ws.Range("C206").Value = 54
$ws.Range("D206").Value = 76182
$ws.Range("C207").Value = 1092
$ws.Range("D207").Value = 1484105
$ws.Range("C208").Value = 1977
$ws.Range("D208").Value = 2455582
$ws.Range("C209").Value = 112
$ws.Range("D209").Value = 136852
$ws.Range("C210").Value = 5848
$ws.Range("D210").Value = 8373954
$ws.Range("C211").Value = 10685
$ws.Range("D211").Value = 14115855
$ws.Range("C212").Value = 1143
$ws.Range("D212").Value = 1556993
$ws.Range("C213").Value = 7158
$ws.Range("D213").Value = 10261834
$ws.Range("C214").Value = 617
$ws.Range("D214").Value = 768027
$ws.Range("C215").Value = 452
$ws.Range("D215").Value = 631700
$ws.Range("C216").Value = 1252
$ws.Range("D216").Value = 1662687
$ws.Range("C217").Value = 3762
$ws.Range("D217").Value = 5095166
$ws.Range("C218").Value = 1803
$ws.Range("D218").Value = 2230747
$ws.Range("C219").Value = 2261
$ws.Range("D219").Value = 2707158
$ws.Range("C220").Value = 4303
$ws.Range("D220").Value = 5732509
$ws.Range("C221").Value = 2206
$ws.Range("D221").Value = 2666304
$ws.Range("C222").Value = 9358
$ws.Range("D222").Value = 11434266
$ws.Range("C223").Value = 3825
$ws.Range("D223").Value = 5280983
$ws.Range("C224").Value = 20
$ws.Range("D224").Value = 28123
$ws.Range("C225").Value = 6315
$ws.Range("D225").Value = 7706173
$ws.Range("C226").Value = 23
$ws.Range("D226").Value = 31259
$ws.Range("C227").Value = 305
$ws.Range("D227").Value = 337345
$ws.Range("C228").Value = 17163
$ws.Range("D228").Value = 24439123
$ws.Range("C229").Value = 24469
$ws.Range("D229").Value = 32465505
$ws.Range("C230").Value = 3108
$ws.Range("D230").Value = 4234738
$ws.Range("C231").Value = 15880
$ws.Range("D231").Value = 22439566
$ws.Range("C232").Value = 1913
$ws.Range("D232").Value = 2468632
$ws.Range("C233").Value = 1032
$ws.Range("D233").Value = 1439001
$ws.Range("C234").Value = 3809
$ws.Range("D234").Value = 5098270
$ws.Range("C235").Value = 11265
$ws.Range("D235").Value = 15108502
$ws.Range("C236").Value = 5879
$ws.Range("D236").Value = 7312793
$ws.Range("C237").Value = 6616
$ws.Range("D237").Value = 7758119
$ws.Range("C238").Value = 14064
$ws.Range("D238").Value = 18608538
$ws.Range("C239").Value = 5523
$ws.Range("D239").Value = 6605209
$ws.Range("C240").Value = 20637
$ws.Range("D240").Value = 25059819
$ws.Range("C241").Value = 3888
$ws.Range("D241").Value = 5231892
$ws.Range("C243").Value = 6421
$ws.Range("D243").Value = 7855898
$ws.Range("C244").Value = 35
$ws.Range("D244").Value = 48661
$ws.Range("C245").Value = 217
$ws.Range("D245").Value = 274076
$ws.Range("C246").Value = 20772
$ws.Range("D246").Value = 29429987
$ws.Range("C247").Value = 26712
$ws.Range("D247").Value = 35105868
$ws.Range("C248").Value = 3350
$ws.Range("D248").Value = 4441402
$ws.Range("C249").Value = 18372
$ws.Range("D249").Value = 25971717
$ws.Range("C250").Value = 2142
$ws.Range("D250").Value = 2732268
$ws.Range("C251").Value = 1073
$ws.Range("D251").Value = 1485582
$ws.Range("C252").Value = 4312
$ws.Range("D252").Value = 5788408
$ws.Range("C253").Value = 12249
$ws.Range("D253").Value = 16609758
$ws.Range("C254").Value = 6176
$ws.Range("D254").Value = 7690308
$ws.Range("C256").Value = 7727
$ws.Range("D256").Value = 9082059
$ws.Range("C257").Value = 13601
$ws.Range("D257").Value = 17941843
$ws.Range("C258").Value = 6169
$ws.Range("D258").Value = 7362860
$ws.Range("C259").Value = 20533
$ws.Range("D259").Value = 23971483
$ws.Range("C260").Value = 1382
$ws.Range("D260").Value = 1911145
$ws.Range("C261").Value = 32
$ws.Range("D261").Value = 38468
$ws.Range("C262").Value = 2693
$ws.Range("D262").Value = 3299119
$ws.Range("C264").Value = 179
$ws.Range("D264").Value = 227177
$ws.Range("C265").Value = 6671
$ws.Range("D265").Value = 9578727
$ws.Range("C266").Value = 11033
$ws.Range("D266").Value = 14562275
$ws.Range("C267").Value = 1404
$ws.Range("D267").Value = 1958174
$ws.Range("C268").Value = 7654
$ws.Range("D268").Value = 11019332
$ws.Range("C269").Value = 1019
$ws.Range("D269").Value = 1296144
$ws.Range("C270").Value = 604
$ws.Range("D270").Value = 852540
$ws.Range("C271").Value = 1944
$ws.Range("D271").Value = 2617388
$ws.Range("C272").Value = 5682
$ws.Range("D272").Value = 7695367
$ws.Range("C273").Value = 2248
$ws.Range("D273").Value = 2881014
$ws.Range("C274").Value = 2827
$ws.Range("D274").Value = 3428015
$ws.Range("C275").Value = 7658
$ws.Range("D275").Value = 10206649
$ws.Range("C276").Value = 2968
$ws.Range("D276").Value = 3583411
$ws.Range("C277").Value = 10692
$ws.Range("D277").Value = 13346867
$ws.Range("C278").Value = 1876
$ws.Range("D278").Value = 2555503
$ws.Range("C279").Value = 6
$ws.Range("D279").Value = 9000
$ws.Range("C280").Value = 6359
$ws.Range("D280").Value = 8182511
$ws.Range("C282").Value = 117
$ws.Range("D282").Value = 149227
$ws.Range("C283").Value = 18204
$ws.Range("D283").Value = 25992346
$ws.Range("C284").Value = 28710
$ws.Range("D284").Value = 38276890
$ws.Range("C285").Value = 6188
$ws.Range("D285").Value = 8643362
$ws.Range("C286").Value = 17945
$ws.Range("D286").Value = 25855515
$ws.Range("C287").Value = 2486
$ws.Range("D287").Value = 3245622
$ws.Range("C288").Value = 1250
$ws.Range("D288").Value = 1766426
$ws.Range("C289").Value = 5256
$ws.Range("D289").Value = 7288417
$ws.Range("C290").Value = 13714
$ws.Range("D290").Value = 18701942
$ws.Range("C291").Value = 7564
$ws.Range("D291").Value = 9743080

# --- Rows 292-313: refreshed counts + NAF section realignment (new section inserted for IDF) ---
$ws.Range("C292").Value = 3
$ws.Range("D292").Value = 4219
$ws.Range("G292").Value = "O"
$ws.Range("H292").Value = "Administration publique"
$ws.Range("C293").Value = 8626
$ws.Range("D293").Value = 10539443
$ws.Range("G293").Value = "P"
$ws.Range("H293").Value = "Enseignement"
$ws.Range("C294").Value = 13403
$ws.Range("D294").Value = 17769263
$ws.Range("G294").Value = "Q"
$ws.Range("H294").Value = "Santé humaine et action sociale"
$ws.Range("C295").Value = 6110
$ws.Range("D295").Value = 7542627
$ws.Range("G295").Value = "R"
$ws.Range("H295").Value = "Arts, spectacles et activités récréatives"
$ws.Range("C296").Value = 19790
$ws.Range("D296").Value = 23651581
$ws.Range("E296").NumberFormat = "@"
$ws.Range("E296").Value = "93"
$ws.Range("F296").Value = "Provence-Alpes-Côte d'Azur"
$ws.Range("G296").Value = "S"
$ws.Range("H296").Value = "Autres activités de services"
$ws.Range("C297").Value = 300
$ws.Range("D297").Value = 411744
$ws.Range("G297").Value = "A"
$ws.Range("H297").Value = "Agriculture, sylviculture et pêche"
$ws.Range("C298").Value = 7152
$ws.Range("D298").Value = 9557847
$ws.Range("G298").Value = "C"
$ws.Range("H298").Value = "Industrie manufacturière"
$ws.Range("C299").Value = 56
$ws.Range("D299").Value = 80598
$ws.Range("G299").Value = "D"
$ws.Range("H299").Value = "Production et distribution d'électricité, de gaz, de vapeur et d'air conditionné"
$ws.Range("C300").Value = 217
$ws.Range("D300").Value = 288251
$ws.Range("G300").Value = "E"
$ws.Range("H300").Value = "Production et distribution d'eau ; assainissement, gestion des déchets et dépollution"
$ws.Range("C301").Value = 21038
$ws.Range("D301").Value = 30414420
$ws.Range("G301").Value = "F"
$ws.Range("H301").Value = "Construction"
$ws.Range("C302").Value = 39211
$ws.Range("D302").Value = 53884041
$ws.Range("G302").Value = "G"
$ws.Range("H302").Value = "Commerce"
$ws.Range("C303").Value = 41622
$ws.Range("D303").Value = 60011477
$ws.Range("G303").Value = "H"
$ws.Range("H303").Value = "Transports et entreposage"
$ws.Range("C304").Value = 24150
$ws.Range("D304").Value = 35561345
$ws.Range("G304").Value = "I"
$ws.Range("H304").Value = "Hébergement et restauration"
$ws.Range("C305").Value = 10480
$ws.Range("D305").Value = 14380751
$ws.Range("G305").Value = "J"
$ws.Range("H305").Value = "Information et communication"
$ws.Range("C306").Value = 2703
$ws.Range("D306").Value = 3875874
$ws.Range("G306").Value = "K"
$ws.Range("H306").Value = "Activités financières et d'assurance"
$ws.Range("C307").Value = 6934
$ws.Range("D307").Value = 9887107
$ws.Range("G307").Value = "L"
$ws.Range("H307").Value = "Activités immobilières"
$ws.Range("C308").Value = 41812
$ws.Range("D308").Value = 57820671
$ws.Range("G308").Value = "M"
$ws.Range("H308").Value = "Activités spécialisées, scientifiques et techniques"
$ws.Range("C309").Value = 11194
$ws.Range("D309").Value = 14992994
$ws.Range("G309").Value = "N"
$ws.Range("H309").Value = "Activités de services administratifs et de soutien"
$ws.Range("C310").Value = 5
$ws.Range("D310").Value = 5675
$ws.Range("G310").Value = "O"
$ws.Range("H310").Value = "Administration publique"
$ws.Range("C311").Value = 12301
$ws.Range("D311").Value = 15303197
$ws.Range("G311").Value = "P"
$ws.Range("H311").Value = "Enseignement"
$ws.Range("C312").Value = 23045
$ws.Range("D312").Value = 31324427
$ws.Range("G312").Value = "Q"
$ws.Range("H312").Value = "Santé humaine et action sociale"
$ws.Range("C313").Value = 15928
$ws.Range("D313").Value = 20541040
$ws.Range("G313").Value = "R"
$ws.Range("H313").Value = "Arts, spectacles et activités récréatives"

# --- New row 314 (additional NAF section for Ile-de-France) ---
$ws.Range("A314").Value = "Fonds de solidarité"
$ws.Range("B314").Value = "VOLET1"
$ws.Range("C314").Value = 25783
$ws.Range("D314").Value = 33526854
$ws.Range("E314").NumberFormat = "@"
$ws.Range("E314").Value = "11"
$ws.Range("F314").Value = "Île-de-France"
$ws.Range("G314").Value = "S"
$ws.Range("H314").Value = "Autres activités de services"
